$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E retain text formatting (many values look numeric,
# e.g. "1.001", "24.397.80", "0.9999" and must not be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '24.397.80'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '1.668.28'
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '313.06'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.3965'
$ws.Range("E7").Value = '  +1.57%  '
$ws.Range("D8").Value = '0.3936'
$ws.Range("E8").Value = '  +2.17%  '
$ws.Range("D9").Value = '52.07'
$ws.Range("E9").Value = '  +6.24%  '
$ws.Range("D10").Value = '1.393'
$ws.Range("E10").Value = '  +3.38%  '
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").Value = '0.08575'
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("E13").Value = '  +2.95%  '
$ws.Range("D14").Value = '7.326'
$ws.Range("E14").Value = '  +2.88%  '
$ws.Range("D15").Value = '7.934'
$ws.Range("E15").Value = '  +5.93%  '
$ws.Range("D16").Value = '0.00001337'
$ws.Range("E16").Value = '  +4.30%  '
$ws.Range("D17").Value = '1.656.07'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").Value = '95.23'
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("D19").Value = '0.07008'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("D21").Value = '6.995'
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("D24").Value = '24.401.10'
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("D25").Value = '2.455'
$ws.Range("E25").Value = '  +5.35%  '
$ws.Range("D26").Value = '3.086'
$ws.Range("E26").Value = '  +14.01%  '
$ws.Range("D27").Value = '22.55'
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").Value = '157.82'
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '142.95'
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").Value = '5.428'
$ws.Range("E30").Value = '  +2.86%  '
$ws.Range("D31").Value = '8.091'
$ws.Range("E31").Value = '  -7.55%  '
$ws.Range("E32").Value = '  +3.70%  '
$ws.Range("D33").Value = '1.843.10'
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("D34").Value = '1.058'
$ws.Range("E34").Value = '  +10.39%  '
$ws.Range("D35").Value = '0.03067'
$ws.Range("E35").Value = '  +5.29%  '
$ws.Range("D36").Value = '0.08253'
$ws.Range("E36").Value = '  +2.98%  '
$ws.Range("D37").Value = '6.899'
$ws.Range("E37").Value = '  -0.96%  '
$ws.Range("E38").Value = '  +12.10%  '
$ws.Range("D39").Value = '0.2770'
$ws.Range("E39").Value = '  +2.66%  '
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.7726'
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '13.81'
$ws.Range("E42").Value = '  +5.59%  '
$ws.Range("D43").Value = '1.446'
$ws.Range("E43").Value = '  -0.93%  '
$ws.Range("D44").Value = '16.67'
$ws.Range("E44").Value = '  +4.15%  '
$ws.Range("D45").Value = '0.7132'
$ws.Range("E45").Value = '  +3.28%  '
$ws.Range("E46").Value = '  +2.51%  '
$ws.Range("D47").Value = '4.138'
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").Value = '0.9999'
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").Value = '0.08434'
$ws.Range("E49").Value = '  +1.06%  '
$ws.Range("D50").Value = '136.78'
$ws.Range("E50").Value = '  +2.28%  '
$ws.Range("D51").Value = '1.267'
$ws.Range("E51").Value = '  +0.36%  '

# Restore the default (Normal) cell style so no stray number-format styling remains.
$ws.Range("D2:E51").Style = "Normal"

